# Add 2022-Q4 data
# 1) Insert a new "2022-Q4" worksheet (as a full copy of the "2022-Q3" sheet so that
#    all styling/number-formats/column widths etc. are preserved) positioned right
#    before "2022-Q3", then overwrite its contents with the new quarter's holdings.
# 2) Insert a new row at the top of the "总计" (summary) sheet with the 2022-Q4 totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: create & populate the new "2022-Q4" sheet
# ---------------------------------------------------------------------------

$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
# NOTE: after Copy(), this same variable now refers to the freshly created
# duplicate sheet (placed right before the original "2022-Q3" sheet).
$q4 = $q3
$q4.Name = "2022-Q4"

# Make sure we have exactly 11 rows (1 header + 10 funds) and 8 columns worth of
# cells available; clear out anything left over from the old (5-row) sheet first.
$q4.Cells.Clear()

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $cell = $q4.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$rows = @(
    @(0, "160642", "鹏华增瑞灵活配置混合（LOF）", "2.22", "85.42", "6.38", "0.1416", 4),
    @(1, "015026", "鹏华增华混合A",               "1.65", "78.57", "3.55", "0.0586", 7),
    @(2, "001675", "江信同福灵活配置混合A",         "0.58", "93.34", "6.90", "0.0400", 7),
    @(3, "620001", "金元顺安宝石动力混合",           "0.46", "56.89", "5.66", "0.0260", 5),
    @(4, "001676", "江信同福灵活配置混合C",         "0.32", "93.34", "6.90", "0.0221", 7),
    @(5, "004927", "中航军民融合精选混合C",         "0.46", "58.03", "2.92", "0.0134", 9),
    @(6, "015027", "鹏华增华混合C",               "0.20", "78.57", "3.55", "0.0071", 7),
    @(7, "004937", "中航混改精选混合C",            "0.07", "74.18", "8.11", "0.0057", 2),
    @(8, "004926", "中航军民融合精选混合A",         "0.13", "58.03", "2.92", "0.0038", 9),
    @(9, "004936", "中航混改精选混合A",            "0.03", "74.18", "8.11", "0.0024", 2)
)

# Pre-format the text-like columns as Text so the numeric-looking strings (fund
# codes, sizes, positions, weights, market values) are preserved verbatim instead
# of being auto-converted into numbers.
$lastRow = 1 + $rows.Count
$q4.Range("B2:B$lastRow").NumberFormat = "@"
$q4.Range("D2:G$lastRow").NumberFormat = "@"

for ($r = 0; $r -lt $rows.Count; $r++) {
    $row = $rows[$r]
    $excelRow = $r + 2
    $q4.Cells.Item($excelRow, 1).Value = $row[0]
    $q4.Cells.Item($excelRow, 2).Value = $row[1]
    $q4.Cells.Item($excelRow, 3).Value = $row[2]
    $q4.Cells.Item($excelRow, 4).Value = $row[3]
    $q4.Cells.Item($excelRow, 5).Value = $row[4]
    $q4.Cells.Item($excelRow, 6).Value = $row[5]
    $q4.Cells.Item($excelRow, 7).Value = $row[6]
    $q4.Cells.Item($excelRow, 8).Value = $row[7]

    $idxCell = $q4.Cells.Item($excelRow, 1)
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1
}

$q4.Range("A1").Select()

# ---------------------------------------------------------------------------
# Part 2: add the 2022-Q4 summary row to "总计"
# ---------------------------------------------------------------------------

$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 10
$total.Range("D2").Value = 0.32

# Column A is just a running 0-based row index; renumber it for every data row
# now that a new row has been inserted at the top of the table.
for ($i = 0; $i -le 5; $i++) {
    $total.Cells.Item($i + 2, 1).Value = $i
}
